$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 - title and link update
$ws.Range("D9").Value = "MBA AI/BigData vs. BSc DS (MSc AI/DS)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/mba-ai-bigdata-vs-bsc-ds/#utm_source=rss&utm_medium=rss&utm_campaign=mba-ai-bigdata-vs-bsc-ds"

# Row 26 - title update only
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 27 - title and link update
$ws.Range("D27").Value = "Apache Beam으로 머신러닝 데이터 파이프라인 구축하기 1편 - 도입과 사용"
$ws.Range("E27").Value = "https://blog.pingpong.us/apache-beam-1/"

# Row 51 - title and link update
$ws.Range("D51").Value = "[세이버메트릭스] MLB 30개 팀의 한글명, 영어명, 팀ID 정리(2022년 기준)"
$ws.Range("E51").Value = "https://bskyvision.com/1278"
